# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Row 20/21 (Uniswap/Litecoin) and row 41/42 (THORChain/MultiversX) swapped rank order,
# so their Coin/Link/Price/Volume cells are rewritten in full; every other row only
# gets refreshed Price (col D) and Volume(1h) (col E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.032.59"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").Value = "2.217.77"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.66"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.34%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").Value = "2.550.83"
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.835"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.45%  "

$ws.Range("D17").Value = "2.200.99"
$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").Value = "41.868.30"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("E19").Value = "  +3.09%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +20.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0798"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.99%  "

$ws.Range("E36").Value = "  -12.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0300"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.39%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.197"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "

$ws.Range("D51").Value = "2.421.12"
$ws.Range("E51").Value = "  -1.56%  "
